$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ------------------------------------------------------------------
# 1) "Patients-" -> bold "Patients" run + plain "-" run
# ------------------------------------------------------------------
$rngPatients = $d.Content
$null = $rngPatients.Find.Execute("Patients", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$rngPatients.Bold = 1

# ------------------------------------------------------------------
# 2) "Doctors-" -> bold "Doctors" run + plain "-" run
# ------------------------------------------------------------------
$rngDoctors = $d.Content
$null = $rngDoctors.Find.Execute("Doctors", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$rngDoctors.Bold = 1

# ------------------------------------------------------------------
# 3) rudransh paragraph: wrap the password "rudransh" in proofErr
#    spellStart/spellEnd and drop the stray _GoBack bookmark (it
#    moves to the new paragraph added in step 4).
# ------------------------------------------------------------------
$pRudransh = $d.Paragraphs.Item(3)
$rngHyperlink3 = $d.Range($pRudransh.Range.Start, $pRudransh.Range.End)
$null = $rngHyperlink3.Find.Execute("rudranshj95@gmail.com", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$afterHyperlink3 = $rngHyperlink3.End
$tailTarget = $d.Range($afterHyperlink3, $pRudransh.Range.End)
$tailTarget.InsertXML(
    "<w:p xmlns:w='$wNs'>" +
    "<w:r><w:br/><w:t>pass-</w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>rudransh</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
)
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 4) New paragraph after the "navdeep" paragraph: hyperlink to
#    jainy@gmail.com, the _GoBack bookmark, "pass-" break-run and a
#    proofErr-wrapped "asdfghjk" run.
# ------------------------------------------------------------------
$pNavdeep = $d.Paragraphs.Item($d.Paragraphs.Count)
$pNavdeep.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item($d.Paragraphs.Count)
$pNew.Range.InsertXML(
    "<w:p xmlns:w='$wNs'>" +
    "<w:r><w:t>jainy@gmail.com</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
    "<w:r><w:br/><w:t>pass-</w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>asdfghjk</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
)

$pNew2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$rngJainy = $d.Range($pNew2.Range.Start, $pNew2.Range.End)
$null = $rngJainy.Find.Execute("jainy@gmail.com", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$null = $d.Hyperlinks.Add($rngJainy, "mailto:jainy@gmail.com", "", "", "jainy@gmail.com")

Write-Host "done"
